$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.364.72"
$ws.Range("D3").Value = "'1.795.17"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'1.005"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'307.23"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "'0.4517"
$ws.Range("D8").Value = "'0.3597"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "'46.34"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").Value = "'0.07083"
$ws.Range("D11").Value = "'0.8847"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "'0.07744"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "'19.47"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'1.810.40"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "'5.285"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "'6.322"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'0.000008511"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "'26.394.60"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'14.24"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "'4.974"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "'2.004.34"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'151.31"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "'17.84"
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").Value = "'2.027"
$ws.Range("E29").Value = "  +4.27%  "
$ws.Range("D30").Value = "'111.87"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").Value = "'4.847"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "'0.08684"
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").Value = "'3.088"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").Value = "'2.761"
$ws.Range("E34").Value = "  +8.91%  "
$ws.Range("D35").Value = "'4.444"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").Value = "'0.7222"
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("D37").Value = "'1.102"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("D38").Value = "'1.004"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "'0.01932"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'0.05085"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("D43").Value = "'0.5063"
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("D44").Value = "'6.827"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "'0.1515"
$ws.Range("E45").Value = "  -4.71%  "
$ws.Range("D46").Value = "'8.014"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "'0.4624"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").Value = "'9.895"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").Value = "'100.61"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "'1.562"
$ws.Range("E51").Value = "  -2.62%  "
